$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.464.31"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").Value = "1.812.15"
$ws.Range("E3").Value = "  +0.67%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.592"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.75%  "

$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("E8").Value = "  +7.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.289"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0675"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0975"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.16%  "

$ws.Range("D12").Value = "2.074.10"
$ws.Range("E12").Value = "  +0.64%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.15%  "

$ws.Range("D14").Value = "1.829.76"
$ws.Range("E14").Value = "  +1.66%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.632"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.74%  "

$ws.Range("D16").Value = "34.452.98"
$ws.Range("E16").Value = "  +0.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.42"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.30"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.92%  "

$ws.Range("D20").Value = "0.0₃0772"
$ws.Range("E20").Value = "  -2.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.95%  "

$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.120"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.72%  "

$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.72%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0517"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.33%  "

$ws.Range("D35").Value = "1.361.47"
$ws.Range("E35").Value = "  -2.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.643"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.22%  "

$ws.Range("E37").Value = "  -0.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.64%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0186"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.65%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.952"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.58%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "81.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.18%  "

$ws.Range("E46").Value = "  +1.32%  "

$ws.Range("D47").Value = "1.975.10"
$ws.Range("E47").Value = "  +0.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.23%  "

$ws.Range("E51").Value = "  -4.20%  "
